$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing column A values (rows 2-3) before they are moved to column B
$profile2 = $ws.Range("A2").Value2
$profile4 = $ws.Range("A3").Value2

# Shift old values into column B
$ws.Range("B2").Value = $profile2
$ws.Range("B3").Value = $profile4

# Write the new values into column A
$ws.Range("A2").Value = "Profile 1"
$ws.Range("A3").Value = "Default"

# Update selection to match the target workbook state
$ws.Range("D8").Select()
